$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (column F) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 7392
$ws1.Range("F11").Value = 7572
$ws1.Range("F15").Value = 6046
$ws1.Range("F24").Value = 276
$ws1.Range("F26").Value = 2086
$ws1.Range("F31").Value = 1048
$ws1.Range("F38").Value = 10

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 12

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 12
$ws4.Range("F14").Value = 7392
$ws4.Range("F16").Value = 7572
$ws4.Range("F18").Value = 6046
$ws4.Range("F26").Value = 276
$ws4.Range("F30").Value = 2086
$ws4.Range("F36").Value = 1048
